# Automatic tracker update
# - Mark rows 166 and 167 as "Fallo" (-1 profit)
# - Append new match rows 175-178

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---- Update existing rows with results ----
Set-TextCell "G166" "Fallo"
$ws.Range("H166").Value = -1

Set-TextCell "G167" "Fallo"
$ws.Range("H167").Value = -1

# ---- Append new rows 175-178 ----

# Row 175
$ws.Range("A175").Value = 14851814
Set-TextCell "B175" "2025-10-17"
Set-TextCell "C175" "Sebastian Korda"
Set-TextCell "D175" "Casper Ruud"
Set-TextCell "E175" "Gana Casper Ruud"
$ws.Range("F175").Value = 1.57

# Row 176
$ws.Range("A176").Value = 14852141
Set-TextCell "B176" "2025-10-16"
Set-TextCell "C176" "Alex Hernandez"
Set-TextCell "D176" "Joao Eduardo Schiessl"
Set-TextCell "E176" "Gana Joao Eduardo Schiessl"
$ws.Range("F176").Value = 1.73

# Row 177
$ws.Range("A177").Value = 14862973
Set-TextCell "B177" "2025-10-16"
Set-TextCell "C177" "Martina Capurro Taborda"
Set-TextCell "D177" "Martina Colmegna"
Set-TextCell "E177" "Gana Martina Colmegna"
$ws.Range("F177").Value = 1.83

# Row 178
$ws.Range("A178").Value = 14893263
Set-TextCell "B178" "2025-10-16"
Set-TextCell "C178" "Bautista Vilicich"
Set-TextCell "D178" "Dante Pagani"
Set-TextCell "E178" "Gana Dante Pagani"
$ws.Range("F178").Value = 3.5
